# Add two new parameter rows to the "Parameters" sheet:
#   row 56: MFTC_WEP_scaling / 1   / description about Winter Energy Payment scaling
#   row 57: WFF_or_Benefit   / Max / description about the WFF-vs-Benefit work decision
#
# These mirror the layout of the existing parameter rows (Parameter | Value | Description
# in columns B | C | D) and reuse the light "section" fill/alignment style used for the
# other top-level parameter rows (e.g. the "Tax/BaseScale" row directly above), minus the
# border lines.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 56: MFTC_WEP_scaling -----------------------------------------
$ws.Range("B56").Value = "MFTC_WEP_scaling"
# Force this to be stored as text ("1"), matching every other Value cell in the
# sheet (they are all text, never numeric), by using a leading quote prefix.
$ws.Range("C56").Value = "'1"
$ws.Range("D56").Value = "How should the Winter Energy Payment be scaled? Average week = 1, Winter week = 12/5, Summer week = 0"

# --- New row 57: WFF_or_Benefit --------------------------------------------
$ws.Range("B57").Value = "WFF_or_Benefit"
$ws.Range("C57").Value = "Max"
$ws.Range("D57").Value = 'What work decision should we assume? Go off-benefit and receive IWTC = "WFF", stay on-benefit = "Benefit", or whichever gives a higher net income = "Max"'

# --- Formatting: match the light-blue "section" style used elsewhere in the
# sheet (same fill as the other Parameter/Value/Description group rows), left
# aligned, no border.
$newRows = $ws.Range("B56:D57")
$newRows.Interior.Color = 14994616
$newRows.HorizontalAlignment = -4131
